$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.155822157859802
$ws.Range("B1").Value = 2.394022703170776
$ws.Range("D1").Value = 2.386180877685547
$ws.Range("E1").Value = 1.225087761878967
